$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "SP0012"
$ws.Range("C2").Value = "SKU0012"
$ws.Range("D2").Value = "Sản phẩm mẫu1"
$ws.Range("H2").Value = 10000
